$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") '44.218.89'
Set-TextCell $ws.Range("E2") '  +1.52%  '

Set-TextCell $ws.Range("D3") '2.245.73'
Set-TextCell $ws.Range("E3") '  +1.14%  '

Set-TextCell $ws.Range("E4") '  +0.15%  '

Set-TextCell $ws.Range("D5") '307.16'
Set-TextCell $ws.Range("E5") '  -1.48%  '

Set-TextCell $ws.Range("D6") '95.91'
Set-TextCell $ws.Range("E6") '  -1.39%  '

Set-TextCell $ws.Range("D7") '0.574'
Set-TextCell $ws.Range("E7") '  +1.42%  '

Set-TextCell $ws.Range("E8") '  +0.23%  '

Set-TextCell $ws.Range("D9") '0.529'
Set-TextCell $ws.Range("E9") '  -0.39%  '

Set-TextCell $ws.Range("D10") '35.28'
Set-TextCell $ws.Range("E10") '  -1.06%  '

Set-TextCell $ws.Range("D11") '0.0817'
Set-TextCell $ws.Range("E11") '  -0.15%  '

Set-TextCell $ws.Range("D12") '7.26'

Set-TextCell $ws.Range("E13") '  +0.36%  '

Set-TextCell $ws.Range("D14") '2.590.64'
Set-TextCell $ws.Range("E14") '  +1.48%  '

Set-TextCell $ws.Range("D15") '2.278.59'
Set-TextCell $ws.Range("E15") '  +2.97%  '

Set-TextCell $ws.Range("D16") '0.835'
Set-TextCell $ws.Range("E16") '  +0.26%  '

Set-TextCell $ws.Range("D17") '13.64'
Set-TextCell $ws.Range("E17") '  -2.65%  '

Set-TextCell $ws.Range("D18") '44.091.36'
Set-TextCell $ws.Range("E18") '  +1.56%  '

Set-TextCell $ws.Range("D19") '0.0₃0972'
Set-TextCell $ws.Range("E19") '  +1.56%  '

Set-TextCell $ws.Range("D20") '6.42'
Set-TextCell $ws.Range("E20") '  +2.49%  '

Set-TextCell $ws.Range("D21") '12.19'
Set-TextCell $ws.Range("E21") '  -5.59%  '

Set-TextCell $ws.Range("D22") '65.58'
Set-TextCell $ws.Range("E22") '  +0.75%  '

Set-TextCell $ws.Range("D23") '237.51'
Set-TextCell $ws.Range("E23") '  +1.15%  '

Set-TextCell $ws.Range("E24") '  +0.70%  '

Set-TextCell $ws.Range("D25") '2.01'
Set-TextCell $ws.Range("E25") '  +0.19%  '

Set-TextCell $ws.Range("E26") '  -0.02%  '

Set-TextCell $ws.Range("D27") '10.02'
Set-TextCell $ws.Range("E27") '  +1.05%  '

Set-TextCell $ws.Range("B28") 'InjectiveProtocol'
Set-TextCell $ws.Range("C28") 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell $ws.Range("D28") '38.42'
Set-TextCell $ws.Range("E28") '  +5.65%  '

Set-TextCell $ws.Range("B29") 'Toncoin'
Set-TextCell $ws.Range("C29") 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell $ws.Range("D29") '2.21'
Set-TextCell $ws.Range("E29") '  +0.99%  '

Set-TextCell $ws.Range("B30") 'Filecoin'
Set-TextCell $ws.Range("C30") 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws.Range("D30") '5.95'
Set-TextCell $ws.Range("E30") '  +0.82%  '

Set-TextCell $ws.Range("B31") 'EthereumClassic'
Set-TextCell $ws.Range("C31") 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell $ws.Range("D31") '20.20'
Set-TextCell $ws.Range("E31") '  +2.16%  '

Set-TextCell $ws.Range("D32") '153.24'
Set-TextCell $ws.Range("E32") '  -3.74%  '

Set-TextCell $ws.Range("D33") '0.0801'
Set-TextCell $ws.Range("E33") '  -2.70%  '

Set-TextCell $ws.Range("E34") '  +3.93%  '

Set-TextCell $ws.Range("E35") '  -1.72%  '

Set-TextCell $ws.Range("E36") '  +3.27%  '

Set-TextCell $ws.Range("E37") '  +1.48%  '

Set-TextCell $ws.Range("E38") '  -5.25%  '

Set-TextCell $ws.Range("D39") '3.51'
Set-TextCell $ws.Range("E39") '  -0.41%  '

Set-TextCell $ws.Range("D40") '14.66'
Set-TextCell $ws.Range("E40") '  -4.12%  '

Set-TextCell $ws.Range("E41") '  -2.99%  '

Set-TextCell $ws.Range("D42") '0.0299'
Set-TextCell $ws.Range("E42") '  -1.84%  '

Set-TextCell $ws.Range("E43") '  +0.23%  '

Set-TextCell $ws.Range("D44") '1.754.04'
Set-TextCell $ws.Range("E44") '  +3.50%  '

Set-TextCell $ws.Range("D45") '83.32'
Set-TextCell $ws.Range("E45") '  +0.01%  '

Set-TextCell $ws.Range("E46") '  +0.13%  '

Set-TextCell $ws.Range("D47") '100.53'
Set-TextCell $ws.Range("E47") '  -0.72%  '

Set-TextCell $ws.Range("D48") '4.97'
Set-TextCell $ws.Range("E48") '  -2.13%  '

Set-TextCell $ws.Range("D49") '8.18'
Set-TextCell $ws.Range("E49") '  +2.64%  '

Set-TextCell $ws.Range("E50") '  -1.80%  '

Set-TextCell $ws.Range("D51") '55.08'
Set-TextCell $ws.Range("E51") '  -1.29%  '
